$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add mapping for electricity: final_energy_carrier id 1 -> primary_energy_carrier id 7
$ws.Range("C2").Value = 7
